# Update Name of Algo
# Applies the value changes described by the diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value  = -12.342
$ws.Range("C10").Value = -12.179
$ws.Range("C12").Value = -12.53
$ws.Range("E13").Value = 12.817
$ws.Range("C18").Value = -12.157
